# Auto-generated Excel COM-interop script to apply market-data refresh
# (scheduled runner update) to the Coeurl_Profits workbook.
#
# For every changed row the columns are:
#   H = currentAveragePrice      I = currentAveragePriceNQ
#   J = currentAveragePriceHQ    K = LevePriceNQ
#   L = LevePriceHQ              M = LeveProfitNQ
#   N = LeveProfitHQ

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 507.4
$ws.Range("J2").Value = 629
$ws.Range("L2").Value = 629
$ws.Range("N2").Value = -855
$ws.Range("H58").Value = 14241.667
$ws.Range("I58").Value = 671.6667
$ws.Range("J58").Value = 21026.666
$ws.Range("K58").Value = 2015.0001
$ws.Range("L58").Value = 63079.99800000001
$ws.Range("M58").Value = -1865.0001
$ws.Range("N58").Value = -63379.99800000001
$ws.Range("H86").Value = 6548.75
$ws.Range("I86").Value = 6500
$ws.Range("J86").Value = 6597.5
$ws.Range("K86").Value = 6500
$ws.Range("L86").Value = 6597.5
$ws.Range("M86").Value = -5377
$ws.Range("N86").Value = -8843.5
$ws.Range("H89").Value = 6548.75
$ws.Range("I89").Value = 6500
$ws.Range("J89").Value = 6597.5
$ws.Range("K89").Value = 32500
$ws.Range("L89").Value = 32987.5
$ws.Range("M89").Value = -26884
$ws.Range("N89").Value = -44219.5
$ws.Range("H98").Value = 2137.5
$ws.Range("I98").Value = 1390.909
$ws.Range("K98").Value = 1390.909
$ws.Range("M98").Value = 107.0909999999999
$ws.Range("H116").Value = 5442.4614
$ws.Range("I116").Value = 5493.2856
$ws.Range("J116").Value = 5383.1665
$ws.Range("K116").Value = 5493.2856
$ws.Range("L116").Value = 5383.1665
$ws.Range("M116").Value = -2051.2856
$ws.Range("N116").Value = -12267.1665
$ws.Range("H122").Value = 2137.5
$ws.Range("I122").Value = 1390.909
$ws.Range("K122").Value = 4172.727000000001
$ws.Range("M122").Value = -1722.727000000001
$ws.Range("H132").Value = 2082.2246
$ws.Range("I132").Value = 2085.4893
$ws.Range("J132").Value = 2005.5
$ws.Range("K132").Value = 6256.467900000001
$ws.Range("L132").Value = 6016.5
$ws.Range("M132").Value = -3726.467900000001
$ws.Range("N132").Value = -11076.5
$ws.Range("H137").Value = 3016.8276
$ws.Range("I137").Value = 3205.2354
$ws.Range("J137").Value = 2749.9167
$ws.Range("K137").Value = 9615.706200000001
$ws.Range("L137").Value = 8249.750100000001
$ws.Range("M137").Value = -7065.706200000001
$ws.Range("N137").Value = -13349.7501
$ws.Range("H138").Value = 5055135
$ws.Range("I138").Value = 3121.625
$ws.Range("J138").Value = 5499268.5
$ws.Range("K138").Value = 9364.875
$ws.Range("L138").Value = 16497805.5
$ws.Range("M138").Value = -4224.875
$ws.Range("N138").Value = -16508085.5

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 779.82355
$ws.Range("J2").Value = 957.1429000000001
$ws.Range("L2").Value = 957.1429000000001
$ws.Range("N2").Value = -1183.1429
$ws.Range("H32").Value = 14050.145
$ws.Range("I32").Value = 7341.2983
$ws.Range("K32").Value = 7341.2983
$ws.Range("M32").Value = -7054.2983
$ws.Range("H45").Value = 10116.417
$ws.Range("I45").Value = 13987.25
$ws.Range("K45").Value = 13987.25
$ws.Range("M45").Value = -13610.25
$ws.Range("H88").Value = 2459.8
$ws.Range("J88").Value = 2766.3333
$ws.Range("L88").Value = 2766.3333
$ws.Range("N88").Value = -3578.3333
$ws.Range("H91").Value = 2459.8
$ws.Range("J91").Value = 2766.3333
$ws.Range("L91").Value = 2766.3333
$ws.Range("N91").Value = -5574.3333
$ws.Range("H116").Value = 779.82355
$ws.Range("J116").Value = 957.1429000000001
$ws.Range("L116").Value = 957.1429000000001
$ws.Range("N116").Value = -5545.1429
$ws.Range("H132").Value = 2352.6304
$ws.Range("I132").Value = 2019.0264
$ws.Range("J132").Value = 3937.25
$ws.Range("K132").Value = 6057.0792
$ws.Range("L132").Value = 11811.75
$ws.Range("M132").Value = -3527.0792
$ws.Range("N132").Value = -16871.75

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 779.82355
$ws.Range("J3").Value = 957.1429000000001
$ws.Range("L3").Value = 957.1429000000001
$ws.Range("N3").Value = -1185.1429
$ws.Range("H94").Value = 3000
$ws.Range("I94").Value = 3000
$ws.Range("K94").Value = 3000
$ws.Range("M94").Value = -2549
$ws.Range("H99").Value = 9952.125
$ws.Range("I99").Value = 6154.25
$ws.Range("K99").Value = 6154.25
$ws.Range("M99").Value = -4656.25
$ws.Range("H107").Value = 2014.0714
$ws.Range("I107").Value = 2065.1538
$ws.Range("J107").Value = 1350
$ws.Range("K107").Value = 2065.1538
$ws.Range("L107").Value = 1350
$ws.Range("M107").Value = -145.1538
$ws.Range("N107").Value = -5190
$ws.Range("H134").Value = 1619.2122
$ws.Range("I134").Value = 1476.069
$ws.Range("J134").Value = 2657
$ws.Range("K134").Value = 4428.207
$ws.Range("L134").Value = 7971
$ws.Range("M134").Value = -1893.207
$ws.Range("N134").Value = -13041

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 27271
$ws.Range("I38").Value = 14500
$ws.Range("J38").Value = 40042
$ws.Range("K38").Value = 14500
$ws.Range("L38").Value = 40042
$ws.Range("M38").Value = -14123
$ws.Range("N38").Value = -40796
$ws.Range("H46").Value = 27271
$ws.Range("I46").Value = 14500
$ws.Range("J46").Value = 40042
$ws.Range("K46").Value = 14500
$ws.Range("L46").Value = 40042
$ws.Range("M46").Value = -14289
$ws.Range("N46").Value = -40464
$ws.Range("H55").Value = 24378.666
$ws.Range("I55").Value = 16527.5
$ws.Range("K55").Value = 16527.5
$ws.Range("M55").Value = -16212.5
$ws.Range("H99").Value = 3111.111
$ws.Range("I99").Value = 3000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -1502
$ws.Range("H126").Value = 3111.111
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530
$ws.Range("H132").Value = 2886.9092
$ws.Range("I132").Value = 2358.75
$ws.Range("J132").Value = 4295.3335
$ws.Range("K132").Value = 7076.25
$ws.Range("L132").Value = 12886.0005
$ws.Range("M132").Value = -4546.25
$ws.Range("N132").Value = -17946.0005
$ws.Range("H134").Value = 11946.115
$ws.Range("I134").Value = 6933.4287
$ws.Range("J134").Value = 32999.4
$ws.Range("K134").Value = 20800.2861
$ws.Range("L134").Value = 98998.20000000001
$ws.Range("M134").Value = -18265.2861
$ws.Range("N134").Value = -104068.2

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 20093.098
$ws.Range("I131").Value = 251544.75
$ws.Range("K131").Value = 754634.25
$ws.Range("M131").Value = -749594.25

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 40001612
$ws.Range("I102").Value = 1722.8096
$ws.Range("K102").Value = 1722.8096
$ws.Range("M102").Value = -100.8096
$ws.Range("H113").Value = 2851.3333
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2851.3333
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = 2851.3333
$ws.Range("N113").Value = -7191.3333
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 2174.739
$ws.Range("J122").Value = 2605.1428
$ws.Range("L122").Value = 7815.428400000001
$ws.Range("N122").Value = -12715.4284
$ws.Range("H126").Value = 15441.556
$ws.Range("I126").Value = 18746.357
$ws.Range("K126").Value = 56239.071
$ws.Range("M126").Value = -53769.071
$ws.Range("H132").Value = 2576.8064
$ws.Range("I132").Value = 1465.6666
$ws.Range("J132").Value = 4115.3076
$ws.Range("K132").Value = 4396.9998
$ws.Range("L132").Value = 12345.9228
$ws.Range("M132").Value = -1866.9998
$ws.Range("N132").Value = -17405.9228

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 22825.938
$ws.Range("I61").Value = 23847.666
$ws.Range("K61").Value = 23847.666
$ws.Range("M61").Value = -23645.666
$ws.Range("H93").Value = 1731.6842
$ws.Range("I93").Value = 1551
$ws.Range("K93").Value = 1551
$ws.Range("M93").Value = -303
$ws.Range("H113").Value = 22825.938
$ws.Range("I113").Value = 23847.666
$ws.Range("K113").Value = 23847.666
$ws.Range("M113").Value = -21677.666

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 14584.143
$ws.Range("I58").Value = 7361.6665
$ws.Range("K58").Value = 7361.6665
$ws.Range("M58").Value = -7053.6665
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("N61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("M61").ClearContents()
$ws.Range("H113").Value = 371790
$ws.Range("I113").Value = 1130.2273
$ws.Range("J113").Value = 2002693
$ws.Range("K113").Value = 3390.6819
$ws.Range("L113").Value = 6008079
$ws.Range("M113").Value = -1220.6819
$ws.Range("N113").Value = -6012419
$ws.Range("H122").Value = 2337.7856
$ws.Range("I122").Value = 2229.1538
$ws.Range("K122").Value = 6687.4614
$ws.Range("M122").Value = -4237.4614
$ws.Range("H126").Value = 4765.5
$ws.Range("J126").Value = 5656.6665
$ws.Range("L126").Value = 16969.9995
$ws.Range("N126").Value = -21909.9995
$ws.Range("H130").Value = 20036.334
$ws.Range("J130").Value = 20036.334
$ws.Range("L130").Value = 20036.334
$ws.Range("N130").Value = -30076.334
$ws.Range("H132").Value = 1756.0149
$ws.Range("I132").Value = 1483.1526
$ws.Range("J132").Value = 3768.375
$ws.Range("K132").Value = 4449.4578
$ws.Range("L132").Value = 11305.125
$ws.Range("M132").Value = -1919.4578
$ws.Range("N132").Value = -16365.125

